# Tabulated the pseudo code
# Indent the statements that sit inside the "if" / "else" branches of the
# begin if / else if / else / end if block so the pseudo code reads as a
# tabulated (indented) block structure.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Statements under the first branch ("homeScore > awayScore ... then")
# get a 3-space indent.
$ws.Range("B31").Value = "   homeScore > awayScore"
$ws.Range("B32").Value = "   then "
$ws.Range("B33").Value = "   winningTeam = A"
$ws.Range("B34").Value = "   teamAPoints = WIN_POINTS"
$ws.Range("B35").Value = "   teamBPoints = 0"

# Statements under the "else if" branch also get a 3-space indent.
$ws.Range("B37").Value = "   homeScore == awayScore"
$ws.Range("B38").Value = "   teamAPoints = DRAW_POINTS"
$ws.Range("B39").Value = "   teamBPoints = DRAW_POINTS"

# Statements under the final "else" branch get a 2-space indent.
$ws.Range("B42").Value = "  winningTeam = B"
$ws.Range("B43").Value = "  teamBPoints = WIN_POINTS"
$ws.Range("B44").Value = "  teamAPoints = 0"
